# Apply the "Actor input / System response" reshuffle on rows 9-11.
#
# Before:
#   Row  9:  C9  blank                        D9  "3. Mostra todos os modelos disponiveis"
#   Row 10:  C10 "4. Seleciona modelo..."      D10 blank
#   Row 11:  C11 blank                        D11 "5. Regista opção"
#
# After:
#   Row  9:  C9  blank                        D9  "3. Verifica modelos disponiveis"
#   Row 10:  C10 blank                        D10 "4. Mostra todos os modelos disponiveis"
#   Row 11:  C11 "5. Seleciona modelo que pretende comprar"   D11 blank

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 - System response column (D) text changes
$ws.Range("D9").Value = "3. Verifica modelos disponiveis"

# Row 10 - move "Mostra todos os modelos disponiveis" text from Actor input (C)
# column to System response (D) column
$ws.Range("C10").Value = $null
$ws.Range("D10").Value = "4. Mostra todos os modelos disponiveis"

# Row 11 - move "Seleciona modelo que pretende comprar" text from System
# response (D) column to Actor input (C) column
$ws.Range("C11").Value = "5. Seleciona modelo que pretende comprar"
$ws.Range("D11").Value = $null

# Update the active selection to match the saved workbook state
$ws.Range("D9").Select()
